$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Travel Request Template (ITIN-2025.1)"

$ws.Range("B3").Value = "<<traveler_name>>"

$ws.Range("B4").Value = "<<business_purpose>>"
$ws.Range("D4").Value = "<<cost_center>>"

$ws.Range("B5").Value = "<<city_state>>"
$ws.Range("F5").Value = "<<destination_zip>>"

$ws.Range("B6").Value = "<<depart_date>>"
$ws.Range("D6").Value = "<<return_date>>"
$ws.Range("F6").Value = "<<event_registration_cost>>"

$ws.Range("B8").Value = "<<flight_pref_outbound.carrier_flight>>"
$ws.Range("C8").Value = "<<flight_pref_outbound.depart_time>>"
$ws.Range("D8").Value = "<<flight_pref_outbound.arrive_time>>"
$ws.Range("E8").Value = "<<flight_pref_outbound.roundtrip_cost>>"

$ws.Range("B9").Value = "<<flight_pref_return.carrier_flight>>"
$ws.Range("C9").Value = "<<flight_pref_return.depart_time>>"
$ws.Range("D9").Value = "<<flight_pref_return.arrive_time>>"
$ws.Range("E9").Value = "<<lowest_cost_roundtrip>>"
$ws.Range("F9").Value = "<<parking_estimate>>"

$ws.Range("B11").Value = "<<hotel.name>>"

$ws.Range("B12").Value = "<<hotel.address>>"
$ws.Range("D12").Value = "<<hotel.city_state>>"
$ws.Range("E12").Value = "<<hotel.nightly_rate>>"
$ws.Range("F12").Value = "<<hotel.nights>>"
$ws.Range("G12").Value = "<<hotel.conference_hotel>>"

$ws.Range("B13").Value = "<<hotel.price_compare_notes>>"

$ws.Range("B14").Value = "<<comparable_hotels[0].name>>"
$ws.Range("C14").Value = "<<comparable_hotels[0].nightly_rate>>"

$ws.Range("B15").Value = "<<ground_transport_pref>>"

$ws.Range("B16").Value = "<<notes>>"

$ws.Range("H20").Formula = "=SUM(H15:H19)"
